# Regenerate the "K" column (G) values for the save_data sheet.
# These values correspond to the recalculated strikeout counts (K) that
# replaced the old "Strike#" derived figures, per the commit message:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 2
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 3
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 2
    20 = 3
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 0
    26 = 3
    27 = 0
    28 = 0
    29 = 0
    30 = 2
    31 = 2
    32 = 2
    33 = 0
    34 = 0
    35 = 3
    36 = 0
    37 = 1
}

foreach ($row in $newK.Keys | Sort-Object) {
    $ws.Range("G$row").Value = $newK[$row]
}
